# Updated cryptos list (price + 1h volume-change refresh, plus a couple of
# rank swaps) per the Nov 28 2024 GitHub Actions data refresh.
#
# Numeric-looking price strings (e.g. "236.96") are written with a leading
# apostrophe so Excel stores them as text (matching the original
# inlineStr/shared-string cells) instead of auto-converting to a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "95.123.91"
$ws.Range("E2").Value = "  -1.21%  "
$ws.Range("D3").Value = "3.577.58"
$ws.Range("E3").Value = "  -0.22%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'236.96"
$ws.Range("E5").Value = "  -1.44%  "
$ws.Range("D6").Value = "'650.37"
$ws.Range("E6").Value = "  +1.98%  "
$ws.Range("E7").Value = "  -0.77%  "
$ws.Range("D8").Value = "'0.399"
$ws.Range("E8").Value = "  -0.62%  "
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").Value = "'1.00"
$ws.Range("E10").Value = "  -2.11%  "
$ws.Range("D11").Value = "3.576.61"
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("D12").Value = "'0.203"
$ws.Range("E12").Value = "  +1.12%  "
$ws.Range("D13").Value = "'42.39"
$ws.Range("E13").Value = "  -1.68%  "
$ws.Range("D14").Value = "'6.50"
$ws.Range("E14").Value = "  +0.72%  "
$ws.Range("D15").Value = "4.242.19"
$ws.Range("E15").Value = "  -0.41%  "
$ws.Range("D16").Value = "94.995.20"
$ws.Range("E16").Value = "  -1.36%  "
$ws.Range("D17").Value = "'0.0000253"
$ws.Range("E17").Value = "  -0.63%  "
$ws.Range("D18").Value = "3.565.19"
$ws.Range("E18").Value = "  -0.51%  "
$ws.Range("D19").Value = "'7.75"
$ws.Range("E19").Value = "  -3.32%  "
$ws.Range("D20").Value = "'12.55"
$ws.Range("E20").Value = "  -4.54%  "
$ws.Range("D21").Value = "'17.89"
$ws.Range("E21").Value = "  -1.55%  "
$ws.Range("D22").Value = "'3.47"
$ws.Range("E22").Value = "  +0.48%  "
$ws.Range("D23").Value = "'507.09"
$ws.Range("E23").Value = "  -1.80%  "
$ws.Range("D24").Value = "'0.477"
$ws.Range("E24").Value = "  -4.40%  "
$ws.Range("D25").Value = "'6.77"
$ws.Range("E25").Value = "  +0.86%  "
$ws.Range("D26").Value = "'0.0000195"
$ws.Range("E26").Value = "  -1.05%  "
$ws.Range("D27").Value = "'95.29"
$ws.Range("E27").Value = "  -1.58%  "
$ws.Range("D28").Value = "'12.50"
$ws.Range("E28").Value = "  +0.70%  "
$ws.Range("D29").Value = "3.769.81"
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "'3.02"
$ws.Range("E30").Value = "  -3.48%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "'0.143"
$ws.Range("E31").Value = "  -0.89%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "'11.48"
$ws.Range("E32").Value = "  -0.83%  "
$ws.Range("E33").Value = "  +0.20%  "
$ws.Range("E34").Value = "  -4.20%  "
$ws.Range("D35").Value = "'0.176"
$ws.Range("E35").Value = "  -2.94%  "
$ws.Range("D36").Value = "'31.76"
$ws.Range("E36").Value = "  +4.83%  "
$ws.Range("D37").Value = "'0.560"
$ws.Range("E37").Value = "  -1.30%  "
$ws.Range("D38").Value = "'1.63"
$ws.Range("E38").Value = "  +8.91%  "
$ws.Range("D39").Value = "'8.50"
$ws.Range("E39").Value = "  +7.17%  "
$ws.Range("D40").Value = "'582.63"
$ws.Range("E40").Value = "  +0.58%  "
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("E42").Value = "  -0.56%  "
$ws.Range("D43").Value = "'0.906"
$ws.Range("D44").Value = "'1.78"
$ws.Range("E44").Value = "  +1.48%  "
$ws.Range("D45").Value = "'2.29"
$ws.Range("E45").Value = "  +4.82%  "
$ws.Range("D46").Value = "'5.70"
$ws.Range("E46").Value = "  +1.09%  "
$ws.Range("D47").Value = "'23.38"
$ws.Range("E47").Value = "  -1.84%  "
$ws.Range("D48").Value = "'33.71"
$ws.Range("E48").Value = "  +31.73%  "
$ws.Range("D49").Value = "'0.0414"
$ws.Range("E49").Value = "  -4.47%  "
$ws.Range("D50").Value = "'3.54"
$ws.Range("E50").Value = "  -0.31%  "
$ws.Range("D51").Value = "'53.27"
$ws.Range("E51").Value = "  -1.15%  "
